$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8089500860585198
$ws.Range("C2").Value = 0.8138528138528138
$ws.Range("D2").Value = 0.8113940440224429
$ws.Range("E2").Value = 1155

$ws.Range("B3").Value = 0.8377952755905512
$ws.Range("C3").Value = 0.8222565687789799
$ws.Range("D3").Value = 0.8299531981279251

$ws.Range("B4").Value = 0.7528373266078184
$ws.Range("C4").Value = 0.758576874205845
$ws.Range("D4").Value = 0.7556962025316456
$ws.Range("E4").Value = 787

$ws.Range("B5").Value = 0.5342857142857143
$ws.Range("C5").Value = 0.5327635327635327
$ws.Range("D5").Value = 0.5335235378031383
$ws.Range("E5").Value = 351

$ws.Range("B6").Value = 0.7673469387755102
$ws.Range("C6").Value = 0.7673469387755102
$ws.Range("D6").Value = 0.7673469387755102
$ws.Range("E6").Value = 0.7673469387755102

$ws.Range("B7").Value = 0.7334671006356509
$ws.Range("C7").Value = 0.7318624474002928
$ws.Range("D7").Value = 0.7326417456212879

$ws.Range("B8").Value = 0.7674857668228967
$ws.Range("C8").Value = 0.7673469387755102
$ws.Range("D8").Value = 0.767394392243536
